# Add three new polling-data rows (Harris Interactive poll, 2021-10-29)
# to the bottom of the PollsData sheet, matching the source commit
# "update w/ harris poll (11/3)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 93 ----
$ws.Range("A93").Value = 27
$ws.Range("B93").Value = 2021
$ws.Range("C93").Value = 9
$ws.Range("D93").Value = 10
$ws.Range("E93").Value = 29
$ws.Range("F93").Value = "harris"
$ws.Range("G93").Value = "online"
$ws.Range("H93").Value = "included"
$ws.Range("I93").Value = 1762
$ws.Range("J93").Value = 1
$ws.Range("K93").Value = 1
$ws.Range("L93").Value = 10
$ws.Range("M93").Value = 2
$ws.Range("N93").Value = 2
$ws.Range("O93").Value = 8
$ws.Range("P93").Value = 5
$ws.Range("Q93").Value = 23
$ws.Range("T93").Value = 14
$ws.Range("U93").Value = "T_0.5"
$ws.Range("V93").Value = 1
$ws.Range("W93").Value = 15
$ws.Range("X93").Value = 17
$ws.Range("Y93").Value = "T_0.5"
$ws.Range("Y93").Font.Color = 0
$ws.Range("AA93").Value = 1
$ws.Range("AB93").Value = "T_0.5"
$ws.Range("AB93").Font.Color = 0

# ---- Row 94 ----
$ws.Range("A94").Value = 27
$ws.Range("B94").Value = 2021
$ws.Range("C94").Value = 9
$ws.Range("D94").Value = 10
$ws.Range("E94").Value = 29
$ws.Range("F94").Value = "harris"
$ws.Range("G94").Value = "online"
$ws.Range("H94").Value = "included"
$ws.Range("I94").Value = 1723
$ws.Range("J94").Value = 1
$ws.Range("K94").Value = 1
$ws.Range("L94").Value = 10
$ws.Range("M94").Value = 2
$ws.Range("N94").Value = 2
$ws.Range("O94").Value = 9
$ws.Range("P94").Value = 5
$ws.Range("Q94").Value = 24
$ws.Range("R94").Value = 10
$ws.Range("U94").Value = "T_0.5"
$ws.Range("V94").Value = 2
$ws.Range("W94").Value = 16
$ws.Range("X94").Value = 17
$ws.Range("Y94").Value = "T_0.5"
$ws.Range("Y94").Font.Color = 0
$ws.Range("AA94").Value = 1
$ws.Range("AB94").Value = "T_0.5"
$ws.Range("AB94").Font.Color = 0

# ---- Row 95 ----
$ws.Range("A95").Value = 27
$ws.Range("B95").Value = 2021
$ws.Range("C95").Value = 9
$ws.Range("D95").Value = 10
$ws.Range("E95").Value = 29
$ws.Range("F95").Value = "harris"
$ws.Range("G95").Value = "online"
$ws.Range("H95").Value = "included"
$ws.Range("I95").Value = 1703
$ws.Range("J95").Value = 1
$ws.Range("K95").Value = 1
$ws.Range("L95").Value = 10
$ws.Range("M95").Value = 2
$ws.Range("N95").Value = 2
$ws.Range("O95").Value = 9
$ws.Range("P95").Value = 5
$ws.Range("Q95").Value = 24
$ws.Range("S95").Value = 9
$ws.Range("U95").Value = "T_0.5"
$ws.Range("V95").Value = 2
$ws.Range("W95").Value = 16
$ws.Range("X95").Value = 18
$ws.Range("Y95").Value = "T_0.5"
$ws.Range("Y95").Font.Color = 0
$ws.Range("AA95").Value = 1
$ws.Range("AB95").Value = "T_0.5"
$ws.Range("AB95").Font.Color = 0

# ---- View state: restore the active selection to where the author left it ----
$ws.Range("Q92").Select()
